# Fix typos in Player architecture diagram
# - Input_Controller -> Input_Controller_Base
# - Movement_Controller -> Movement_Controller_Base (split into 3 runs, as PowerPoint
#   does when the run is retyped / spell-checked)
# - refresh the auto date placeholders (datetimeFigureOut) on every layout + the master

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide content fixes
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange

    if ($shp.Name -eq "Rectangle 6") {
        # "Input_Controller" -> "Input_Controller_Base"
        $full = $tr.Text
        if ($full.IndexOf("Input_Controller") -eq 0) {
            $sub = $tr.Characters(1, 17)
            $sub.Text = "Input_Controller_Base"
        }
    }

    if ($shp.Name -eq "Rectangle 7") {
        # "Movement_Controller" -> "Movement" + "_" + "Controller_Base"
        $full = $tr.Text
        if ($full.IndexOf("Movement_Controller") -eq 0) {
            $c3 = $tr.Characters(10, 10)
            $c3.Text = "Controller_Base"
            $c2 = $tr.Characters(9, 1)
            $c2.Text = "_"
            $c1 = $tr.Characters(1, 8)
            $c1.Text = "Movement"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Bump the auto-update date placeholder text on every slide layout + the
#    slide master (mirrors what PowerPoint re-caches on save).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name.IndexOf("Date Placeholder") -eq 0) {
            if ($shp.HasTextFrame) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq "7/22/2025") {
                    $tr.Text = "7/23/2025"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i)
}
